$d = $word.ActiveDocument
$sec = $d.Sections(1)
$ftr = $sec.Footers.Item(1)

$rng2 = $ftr.Range
$rng2.Find.Execute("Spring 2015", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
Write-Host "after delete: [$($ftr.Range.Text)]"

# Try re-finding tab and using Characters/Words collection item-level insert
$tab = [char]9
$rng3 = $ftr.Range
$rng3.Find.Execute($tab) | Out-Null
Write-Host "tab range start end:" $rng3.Start $rng3.End

$rng4 = $ftr.Range
$rng4.Start = $rng3.End
$rng4.End = $rng3.End
Write-Host "rng4 start end (after manual set):" $rng4.Start $rng4.End
$rng4.Text = "Spring 2016"
Write-Host "footer text after set .Text on rng4: [$($ftr.Range.Text)]"
